$wb = $excel.ActiveWorkbook

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 15238.814
$ws.Range("I32").Value = 9858.339
$ws.Range("J32").Value = 56937.5
$ws.Range("K32").Value = 9858.339
$ws.Range("L32").Value = 56937.5
$ws.Range("M32").Value = -9571.339
$ws.Range("N32").Value = -57511.5
# Row 88: The Mast Chance
$ws.Range("H88").Value = 2019.1904
$ws.Range("I88").Value = 2022.0625
$ws.Range("J88").Value = 2010
$ws.Range("K88").Value = 2022.0625
$ws.Range("L88").Value = 2010
$ws.Range("M88").Value = -1616.0625
$ws.Range("N88").Value = -2822
# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 2019.1904
$ws.Range("I91").Value = 2022.0625
$ws.Range("J91").Value = 2010
$ws.Range("K91").Value = 2022.0625
$ws.Range("L91").Value = 2010
$ws.Range("M91").Value = -618.0625
$ws.Range("N91").Value = -4818

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 6307.5
$ws.Range("I86").Value = 5414.357
$ws.Range("J86").Value = 7870.5
$ws.Range("K86").Value = 5414.357
$ws.Range("L86").Value = 7870.5
$ws.Range("M86").Value = -4291.357
$ws.Range("N86").Value = -10116.5
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 6307.5
$ws.Range("I89").Value = 5414.357
$ws.Range("J89").Value = 7870.5
$ws.Range("K89").Value = 27071.785
$ws.Range("L89").Value = 39352.5
$ws.Range("M89").Value = -21455.785
$ws.Range("N89").Value = -50584.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3771.7932
$ws.Range("I31").Value = 3261.0476
$ws.Range("J31").Value = 5112.5
$ws.Range("K31").Value = 3261.0476
$ws.Range("L31").Value = 5112.5
$ws.Range("M31").Value = -2966.0476
$ws.Range("N31").Value = -5702.5
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3771.7932
$ws.Range("I34").Value = 3261.0476
$ws.Range("J34").Value = 5112.5
$ws.Range("K34").Value = 3261.0476
$ws.Range("L34").Value = 5112.5
$ws.Range("M34").Value = -3059.0476
$ws.Range("N34").Value = -5516.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 712.375
$ws.Range("I5").Value = 436.89474
$ws.Range("J5").Value = 961.619
$ws.Range("K5").Value = 1310.68422
$ws.Range("L5").Value = 2884.857
$ws.Range("M5").Value = -1198.68422
$ws.Range("N5").Value = -3108.857
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 497.85294
$ws.Range("I113").Value = 517.08
$ws.Range("J113").Value = 444.44446
$ws.Range("K113").Value = 1551.24
$ws.Range("L113").Value = 1333.33338
$ws.Range("M113").Value = 618.7599999999998
$ws.Range("N113").Value = -5673.33338
# Row 126: Imperial Palate
$ws.Range("H126").Value = 9336.25
$ws.Range("I126").Value = 10606.154
$ws.Range("J126").Value = 3833.3333
$ws.Range("K126").Value = 31818.462
$ws.Range("L126").Value = 11499.9999
$ws.Range("M126").Value = -26878.462
$ws.Range("N126").Value = -21379.9999
# Row 130: Blast from the Pasta
$ws.Range("H130").Value = 500
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 4156.1113
$ws.Range("I131").Value = 3507.7
$ws.Range("J131").Value = 4966.625
$ws.Range("K131").Value = 10523.1
$ws.Range("L131").Value = 14899.875
$ws.Range("M131").Value = -5483.099999999999
$ws.Range("N131").Value = -24979.875
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 712.375
$ws.Range("I135").Value = 436.89474
$ws.Range("J135").Value = 961.619
$ws.Range("K135").Value = 3932.05266
$ws.Range("L135").Value = 8654.571
$ws.Range("M135").Value = -1397.05266
$ws.Range("N135").Value = -13724.571
# Row 139: Najoothie
$ws.Range("H139").Value = 1475.3846
$ws.Range("I139").Value = 940
$ws.Range("K139").Value = 2820
$ws.Range("M139").Value = 2320

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1259.4546
$ws.Range("I97").Value = 988.05884
$ws.Range("J97").Value = 2182.2
$ws.Range("K97").Value = 988.05884
$ws.Range("L97").Value = 2182.2
$ws.Range("M97").Value = -492.05884
$ws.Range("N97").Value = -3174.2
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 4754.107
$ws.Range("I102").Value = 5278.1816
$ws.Range("J102").Value = 2832.5
$ws.Range("K102").Value = 5278.1816
$ws.Range("L102").Value = 2832.5
$ws.Range("M102").Value = -3656.1816
$ws.Range("N102").Value = -6076.5
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 6421.722
$ws.Range("I113").Value = 9752.909
$ws.Range("J113").Value = 1187
$ws.Range("K113").Value = 9752.909
$ws.Range("L113").Value = 1187
$ws.Range("M113").Value = -7582.909
$ws.Range("N113").Value = -5527

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 510.76923
$ws.Range("I22").Value = 410
$ws.Range("J22").Value = 541
$ws.Range("K22").Value = 410
$ws.Range("L22").Value = 541
$ws.Range("M22").Value = -115
$ws.Range("N22").Value = -1131
# Row 27: Fire and Hide
$ws.Range("H27").Value = 510.76923
$ws.Range("I27").Value = 410
$ws.Range("J27").Value = 541
$ws.Range("K27").Value = 410
$ws.Range("L27").Value = 541
$ws.Range("M27").Value = -303
$ws.Range("N27").Value = -755

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 43: Walk Softly and Carry a Big Halberd
$ws.Range("H43").Value = 14580
$ws.Range("J43").Value = 14580
$ws.Range("L43").Value = 14580
$ws.Range("N43").Value = -14878
# Row 52: Party Animals
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
# Row 76: Finger on the Pulse
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
# Row 79: Chirurgeon Hand in Glove (L)
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 2733.3333
$ws.Range("I81").Value = 1416.6666
$ws.Range("J81").Value = 5366.6665
$ws.Range("K81").Value = 2833.3332
$ws.Range("L81").Value = 10733.333
$ws.Range("M81").Value = -1772.3332
$ws.Range("N81").Value = -12855.333
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 2733.3333
$ws.Range("I84").Value = 1416.6666
$ws.Range("J84").Value = 5366.6665
$ws.Range("K84").Value = 14166.666
$ws.Range("L84").Value = 53666.665
$ws.Range("M84").Value = -8862.666000000001
$ws.Range("N84").Value = -64274.665
# Row 103: To the Tops
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("N103").ClearContents()
# Row 109: Turban in Training
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").ClearContents()
# Row 113: A Tender Table
$ws.Range("H113").Value = 645.2857
$ws.Range("I113").Value = 627.7143
$ws.Range("J113").Value = 662.8570999999999
$ws.Range("K113").Value = 1883.1429
$ws.Range("L113").Value = 1988.5713
$ws.Range("M113").Value = 286.8571000000002
$ws.Range("N113").Value = -6328.5713
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2900.889
$ws.Range("I122").Value = 2978.1538
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 8934.4614
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -6484.4614
$ws.Range("N122").Value = -13000
# Row 133: Begin with the Basics
$ws.Range("H133").Value = 27666.666
$ws.Range("J133").Value = 27666.666
$ws.Range("L133").Value = 27666.666
$ws.Range("N133").Value = -37786.666
# Row 137: Traditional Trousers
$ws.Range("H137").Value = 34380.625
$ws.Range("I137").Value = 32325
$ws.Range("J137").Value = 35065.832
$ws.Range("K137").Value = 32325
$ws.Range("L137").Value = 35065.832
$ws.Range("M137").Value = -27225
$ws.Range("N137").Value = -45265.832
# Row 139: Cruel Climates
$ws.Range("H139").Value = 43785.832
$ws.Range("J139").Value = 43785.832
$ws.Range("L139").Value = 43785.832
$ws.Range("N139").Value = -54065.832
# Row 141: Silk for Sunperch
$ws.Range("J141").Value = 45678.75
$ws.Range("L141").Value = 45678.75
$ws.Range("N141").Value = -56038.75
